$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, pushing existing rows 5-8 down to 6-9
$ws.Rows.Item(5).Insert()

# Columns A-L are identical across all the existing entries for this
# market/product, so copy them from the row right below (now row 6,
# the former row 5) into the newly inserted row 5.
for ($c = 1; $c -le 12; $c++) {
    $ws.Cells.Item(5, $c).Value2 = $ws.Cells.Item(6, $c).Value2
}

# Apply the same number format as the other date cells in column D to the new D5 cell
$ws.Cells.Item(5, 4).NumberFormat = $ws.Cells.Item(6, 4).NumberFormat

# Set the specific new values for row 5
$ws.Cells.Item(5, 4).Value2 = 44880
$ws.Cells.Item(5, 13).Value2 = 200
$ws.Cells.Item(5, 14).Value2 = 33000
$ws.Cells.Item(5, 15).Value2 = 34000
$ws.Cells.Item(5, 16).Value2 = 33500
$ws.Cells.Item(5, 17).Value2 = "$/caja 10 kilos"
$ws.Cells.Item(5, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(5, 19).Value2 = 3350
$ws.Cells.Item(5, 20).Value2 = 10
